$d = $word.ActiveDocument

function GetRange($needle) {
    $r = $d.Content.Duplicate
    $r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return $r
}

# ---------------------------------------------------------------------------
# Hunk 1: "No." + " " + "file no." (3 runs, first two bold Arial) collapse
# into a single plain run "No. file no." (formatting of the trailing run,
# which had no explicit rPr).
# ---------------------------------------------------------------------------
$full1 = GetRange("No. file no.")
$s1 = $full1.Start
$e1 = $full1.End
$ins1 = $d.Range($e1, $e1)
$ins1.InsertAfter("No. file no.")
$d.Range($s1, $e1).Delete()

# ---------------------------------------------------------------------------
# Hunk 2: "  " + "31-07-2025" (two plain runs) -> single plain run
# "  01-08-2025" (date bumped by one day, leading spaces preserved).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("  31-07-2025", $true, $false, $false, $false, $false, `
    $true, 1, $false, "  01-08-2025", 2) | Out-Null

# ---------------------------------------------------------------------------
# Hunk 3: "Firm X" (plain) + "," (bold Arial) -> single plain run "Firm X,"
# (formatting of the leading run, which had no explicit rPr).
# ---------------------------------------------------------------------------
$full3 = GetRange("Firm X,")
$s3 = $full3.Start
$e3 = $full3.End
$ins3 = $d.Range($s3, $s3)
$ins3.InsertBefore("Firm X,")
$newLen3 = "Firm X,".Length
$d.Range($s3 + $newLen3, $e3 + $newLen3).Delete()

# ---------------------------------------------------------------------------
# Hunk 4: "firm address" -> "<<FIRM_ADDRESS>>" (same run, text only).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("firm address", $true, $false, $false, $false, $false, `
    $true, 1, $false, "<<FIRM_ADDRESS>>", 2) | Out-Null

# ---------------------------------------------------------------------------
# Hunk 5: 5 runs forming "This office's LOA No. lOA NO. dated 01-07-2025."
# collapse into a single plain run
# "This office's LOA No. {{LOA_NO}} dated {{LOA_DATE}}."
# (formatting of the "lOA NO." run, which had no explicit rPr).
# ---------------------------------------------------------------------------
$full5 = GetRange("This office's LOA No. lOA NO. dated 01-07-2025.")
$fullStart5 = $full5.Start
$fullEnd5 = $full5.End

$anchor5 = GetRange("lOA NO.")
$anchorEnd5 = $anchor5.End

$newText5 = "This office's LOA No. {{LOA_NO}} dated {{LOA_DATE}}."
$insPoint5 = $d.Range($anchorEnd5, $anchorEnd5)
$insPoint5.InsertBefore($newText5)
$newLen5 = $newText5.Length

$tailStart5 = $anchorEnd5 + $newLen5
$tailEnd5 = $fullEnd5 + $newLen5
$d.Range($tailStart5, $tailEnd5).Delete()

$d.Range($fullStart5, $anchorEnd5).Delete()
